$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: Coin name / Link / Price / Volume(1h) columns refreshed from the
# latest coinranking.com snapshot. Rows 28/29 and 45/46 swap contents (rank reorder).
$rowUpdates = @{
  2 = @{ D='27.559.77'; E='  -1.51%  ' }
  3 = @{ D='1.666.16'; E='  -3.50%  ' }
  4 = @{ E='  +0.11%  ' }
  5 = @{ D='215.18'; E='  -1.56%  ' }
  6 = @{ E='  -1.94%  ' }
  7 = @{ E='  +0.11%  ' }
  8 = @{ D='23.57'; E='  -2.18%  ' }
  9 = @{ D='0.263'; E='  -0.91%  ' }
  10 = @{ E='  -2.11%  ' }
  11 = @{ D='0.0882'; E='  -2.12%  ' }
  12 = @{ D='1.902.75'; E='  -3.44%  ' }
  13 = @{ D='1.671.36'; E='  -3.16%  ' }
  14 = @{ E='  -2.74%  ' }
  15 = @{ D='0.555'; E='  -2.20%  ' }
  16 = @{ D='66.25'; E='  -2.46%  ' }
  17 = @{ D='250.17'; E='  +2.48%  ' }
  18 = @{ D='27.590.77'; E='  -1.27%  ' }
  19 = @{ E='  -3.28%  ' }
  20 = @{ E='  -4.20%  ' }
  21 = @{ D='1.00'; E='  +0.14%  ' }
  22 = @{ E='  -3.28%  ' }
  23 = @{ E='  -4.82%  ' }
  24 = @{ E='  -5.53%  ' }
  25 = @{ D='146.64'; E='  -1.91%  ' }
  26 = @{ D='16.50'; E='  -2.05%  ' }
  27 = @{ E='  -5.09%  ' }
  28 = @{ B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.00'; E='  +0.17%  ' }
  29 = @{ B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.112'; E='  -2.38%  ' }
  30 = @{ E='  +4.13%  ' }
  31 = @{ D='0.0507'; E='  -0.95%  ' }
  32 = @{ E='  -2.88%  ' }
  33 = @{ D='1.474.31'; E='  -1.08%  ' }
  34 = @{ E='  -5.66%  ' }
  35 = @{ D='1.58'; E='  -5.18%  ' }
  36 = @{ D='0.940'; E='  -2.15%  ' }
  37 = @{ E='  -1.16%  ' }
  38 = @{ D='0.575'; E='  -6.28%  ' }
  39 = @{ E='  -2.40%  ' }
  40 = @{ D='69.68'; E='  -2.71%  ' }
  41 = @{ D='1.02'; E='  -4.46%  ' }
  42 = @{ D='1.00'; E='  +0.06%  ' }
  43 = @{ D='5.42'; E='  -7.27%  ' }
  44 = @{ D='1.810.73'; E='  -3.40%  ' }
  45 = @{ B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.21'; E='  -3.50%  ' }
  46 = @{ B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.788'; E='  -0.87%  ' }
  47 = @{ E='  -2.03%  ' }
  48 = @{ D='89.35'; E='  -2.16%  ' }
  49 = @{ E='  -2.13%  ' }
  50 = @{ D='41.93'; E='  +15.57%  ' }
  51 = @{ E='  -3.32%  ' }
}

foreach ($row in $rowUpdates.Keys) {
  $fields = $rowUpdates[$row]
  foreach ($col in $fields.Keys) {
    $value = $fields[$col]
    $cell = $ws.Range("$col$row")
    if ($col -eq "D" -and $value -match "^-?[0-9]+(\.[0-9]+)?$") {
      # Force text so Excel does not auto-convert numeric-looking price strings
      # (e.g. '1.00', '215.18') into real numbers, matching the original inline-string cells.
      $cell.Value = "'" + $value
    } else {
      $cell.Value = $value
    }
  }
}
